$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(59, "2023-12-07 11:02:41", 0.0004),
    @(60, "2023-12-07 11:02:55", 0.0006000000000000001),
    @(61, "2023-12-07 11:03:30", 0.0012),
    @(62, "2023-12-07 11:03:36", 0.0004),
    @(63, "2023-12-07 11:03:56", 0.0008)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
}
